# Adds a new "PRESUPUESTO" (G) column to the "VENTA MENSUAL" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# 1) Copy the formatting of column F (currency data + bold header + right
#    aligned totals row) into the new column G, so the new cells inherit the
#    same styles already used by the sheet (header style, currency style,
#    totals style).
$ws.Range("F1:F264").Copy() | Out-Null
$ws.Range("G1:G264").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) Column G should be 17 units wide (raw OOXML width). Excel's
#    ColumnWidth COM property adds a constant ~0.8333 padding on top of the
#    raw stored width, so compensate for that offset.
$ws.Columns.Item(7).ColumnWidth = 17 - 0.8333333333333334

# 3) Header label.
$ws.Cells.Item(1, 7).Value = "PRESUPUESTO"

# 4) Budget ("PRESUPUESTO") values for every sales row (2-263). Most rows
#    are 0; a handful of advisors have an assigned monthly budget.
$gValues = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1000,0,3000,5000,0,6000,0,7000,0,6000,6000,1000,400,6500,0,4000,0,500,4000,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$rowCount = $gValues.Length
$data = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i, 0] = $gValues[$i]
}
$ws.Range("G2:G263").Value = $data

# 5) Totals row (sum of the PRESUPUESTO column).
$ws.Cells.Item(264, 7).Value = 50400
